$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.002.66"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.919.57"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4578"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3807"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07743"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9763"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "1.923.42"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.703"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.937"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07019"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009480"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "29.015.40"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.342"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "2.152.69"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.062"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.594"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.826"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09316"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8550"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.077"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.237"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.021"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05672"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.102"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.406"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5475"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1750"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.332"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002855"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.171"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5164"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06914"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.751"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.63%  "
